# Update of Slovakia Super Liga base (06-04-2024 15:39)
# - Rows 15-17: the three fixtures played on 2023-08-12 were re-sequenced
#   (ids kept in place, underlying match data rotates up by one row).
# - Rows 20-21, 59-60, 124-125: pairs of fixtures swap their match data.
# - Row 153: the still-unplayed "MFK Ruzomberok vs FC Spartak Trnava" entry
#   is replaced by an updated "Dunajska Streda vs FK Zeleziarne Podbrezova"
#   fixture with refreshed odds.
# - Rows 154-158 (the other not-yet-played placeholder fixtures) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 15, 16, 17: rotate the B:AC payload up by one (17 -> 15 -> 16 -> 17) ---
$r15 = $ws.Range("B15:AC15")
$r16 = $ws.Range("B16:AC16")
$r17 = $ws.Range("B17:AC17")
$v15 = $r15.Value()
$v16 = $r16.Value()
$v17 = $r17.Value()
$r15.Value = $v17
$r16.Value = $v15
$r17.Value = $v16

# --- Rows 20 <-> 21 ---
$r20 = $ws.Range("B20:AC20")
$r21 = $ws.Range("B21:AC21")
$v20 = $r20.Value()
$v21 = $r21.Value()
$r20.Value = $v21
$r21.Value = $v20

# --- Rows 59 <-> 60 ---
$r59 = $ws.Range("B59:AC59")
$r60 = $ws.Range("B60:AC60")
$v59 = $r59.Value()
$v60 = $r60.Value()
$r59.Value = $v60
$r60.Value = $v59

# --- Rows 124 <-> 125 ---
$r124 = $ws.Range("B124:AC124")
$r125 = $ws.Range("B125:AC125")
$v124 = $r124.Value()
$v125 = $r125.Value()
$r124.Value = $v125
$r125.Value = $v124

# --- Row 153: replace with the refreshed Dunajska Streda vs FK Zeleziarne Podbrezova line ---
$ws.Cells.Item(153, 2).Value = 7958336
$ws.Cells.Item(153, 5).Value = 45389.52083333334
$ws.Cells.Item(153, 6).Value = "Dunajska Streda"
$ws.Cells.Item(153, 7).Value = "FK Zeleziarne Podbrezova"
$ws.Cells.Item(153, 11).Value = 1.833
$ws.Cells.Item(153, 12).Value = 3.6
$ws.Cells.Item(153, 13).Value = 3.75
$ws.Cells.Item(153, 14).Value = 1.95
$ws.Cells.Item(153, 15).Value = 3.5
$ws.Cells.Item(153, 16).Value = 4
$ws.Cells.Item(153, 17).Value = -0.5
$ws.Cells.Item(153, 18).Value = 1.975
$ws.Cells.Item(153, 19).Value = 1.825
$ws.Cells.Item(153, 20).Value = 2.5
$ws.Cells.Item(153, 21).Value = 1.8
$ws.Cells.Item(153, 22).Value = 2
$ws.Cells.Item(153, 23).Value = 0
$ws.Cells.Item(153, 24).Value = 0
$ws.Cells.Item(153, 25).Value = 0
$ws.Cells.Item(153, 26).Value = 0
$ws.Cells.Item(153, 27).Value = 0

# --- Remove the now-superseded placeholder rows 154-158 ---
$ws.Range("A154:A158").EntireRow.Delete()
